$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 80
$ws.Range("H80").Value = 610.7857
$ws.Range("I80").Value = 608
$ws.Range("J80").Value = 614.5
$ws.Range("K80").Value = 1824
$ws.Range("L80").Value = 1843.5
$ws.Range("M80").Value = -826
$ws.Range("N80").Value = -3839.5
# Row 83
$ws.Range("H83").Value = 610.7857
$ws.Range("I83").Value = 608
$ws.Range("J83").Value = 614.5
$ws.Range("K83").Value = 5472
$ws.Range("L83").Value = 5530.5
$ws.Range("M83").Value = -480
$ws.Range("N83").Value = -15514.5
# Row 86
$ws.Range("H86").Value = 4266.6665
$ws.Range("I86").Value = 3900
$ws.Range("K86").Value = 3900
$ws.Range("M86").Value = -2777
# Row 89
$ws.Range("H89").Value = 4266.6665
$ws.Range("I89").Value = 3900
$ws.Range("K89").Value = 19500
$ws.Range("M89").Value = -13884
# Row 132
$ws.Range("H132").Value = 12350.833
$ws.Range("I132").Value = 11451.5
$ws.Range("K132").Value = 34354.5
$ws.Range("M132").Value = -31824.5
# Row 138
$ws.Range("H138").Value = 2916.1333
$ws.Range("I138").Value = 1220.8889
$ws.Range("J138").Value = 5459
$ws.Range("K138").Value = 3662.6667
$ws.Range("L138").Value = 16377
$ws.Range("M138").Value = 1477.3333
$ws.Range("N138").Value = -26657
# Row 141
$ws.Range("H141").Value = 833.96295
$ws.Range("I141").Value = 833.96295
$ws.Range("K141").Value = 2501.88885
$ws.Range("M141").Value = 2678.11115

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 61
$ws.Range("H61").Value = 2489
$ws.Range("I61").Value = 2314.2727
$ws.Range("K61").Value = 2314.2727
$ws.Range("M61").Value = -2102.2727
# Row 122
$ws.Range("H122").Value = 2469.7778
$ws.Range("I122").Value = 1950
$ws.Range("J122").Value = 3509.3333
$ws.Range("K122").Value = 5850
$ws.Range("L122").Value = 10527.9999
$ws.Range("M122").Value = -3400
$ws.Range("N122").Value = -15427.9999
# Row 130
$ws.Range("H130").Value = 48990
$ws.Range("J130").Value = 48990
$ws.Range("L130").Value = 48990
$ws.Range("N130").Value = -59030
# Row 136
$ws.Range("H136").Value = 2489
$ws.Range("I136").Value = 2314.2727
$ws.Range("K136").Value = 6942.8181
$ws.Range("M136").Value = -4392.8181

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 86
$ws.Range("H86").Value = 2243.111
$ws.Range("J86").Value = 1874.75
$ws.Range("L86").Value = 1874.75
$ws.Range("N86").Value = -4120.75
# Row 89
$ws.Range("H89").Value = 2243.111
$ws.Range("J89").Value = 1874.75
$ws.Range("L89").Value = 9373.75
$ws.Range("N89").Value = -20605.75
# Row 99
$ws.Range("H99").Value = 1344.9
$ws.Range("I99").Value = 1362.1111
$ws.Range("K99").Value = 1362.1111
$ws.Range("M99").Value = 135.8888999999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 22
$ws.Range("H22").Value = 669
$ws.Range("I22").Value = 659.8461
$ws.Range("K22").Value = 659.8461
$ws.Range("M22").Value = -309.8461
# Row 44
$ws.Range("H44").Value = 14999.333
$ws.Range("I44").Value = 5000
$ws.Range("K44").Value = 5000
$ws.Range("M44").Value = -4558
# Row 58
$ws.Range("H58").Value = 2439.2068
$ws.Range("I58").Value = 1430.28
$ws.Range("K58").Value = 1430.28
$ws.Range("M58").Value = -1227.28
# Row 132
$ws.Range("H132").Value = 2004.5682
$ws.Range("I132").Value = 2055.195
$ws.Range("K132").Value = 6165.585000000001
$ws.Range("M132").Value = -3635.585000000001
# Row 136
$ws.Range("H136").Value = 2439.2068
$ws.Range("I136").Value = 1430.28
$ws.Range("K136").Value = 4290.84
$ws.Range("M136").Value = -1740.84

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 14
$ws.Range("H14").Value = 441.5
$ws.Range("I14").Value = 441.5
$ws.Range("K14").Value = 1324.5
$ws.Range("M14").Value = -1151.5
# Row 34
$ws.Range("H34").Value = 2185
$ws.Range("J34").Value = 3582.8333
$ws.Range("L34").Value = 10748.4999
$ws.Range("N34").Value = -10916.4999
# Row 75
$ws.Range("H75").Value = 945.6
$ws.Range("J75").Value = 1003.75
$ws.Range("L75").Value = 3011.25
$ws.Range("N75").Value = -5007.25
# Row 78
$ws.Range("H78").Value = 945.6
$ws.Range("J78").Value = 1003.75
$ws.Range("L78").Value = 9033.75
$ws.Range("N78").Value = -19017.75
# Row 129
$ws.Range("H129").Value = 1315.8182
$ws.Range("I129").Value = 684.875
$ws.Range("K129").Value = 2054.625
$ws.Range("M129").Value = 2945.375
# Row 137
$ws.Range("H137").Value = 4250
$ws.Range("J137").Value = 4250
$ws.Range("L137").Value = 12750
$ws.Range("N137").Value = -22950

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 122
$ws.Range("H122").Value = 2592.3333
$ws.Range("I122").Value = 1401.75
$ws.Range("J122").Value = 3544.8
$ws.Range("K122").Value = 4205.25
$ws.Range("L122").Value = 10634.4
$ws.Range("M122").Value = -1755.25
$ws.Range("N122").Value = -15534.4
# Row 132
$ws.Range("H132").Value = 21037.574
$ws.Range("I132").Value = 23662.447
$ws.Range("K132").Value = 70987.341
$ws.Range("M132").Value = -68457.341

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 122
$ws.Range("H122").Value = 3000
$ws.Range("I122").Value = 3000
$ws.Range("K122").Value = 9000
$ws.Range("M122").Value = -6550
# Row 132
$ws.Range("H132").Value = 3219
$ws.Range("I132").Value = 1714.2307
$ws.Range("K132").Value = 5142.6921
$ws.Range("M132").Value = -2612.6921
# Row 136
$ws.Range("H136").Value = 2999.6
$ws.Range("I136").Value = 2999.5
$ws.Range("K136").Value = 8998.5
$ws.Range("M136").Value = -6448.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 24
$ws.Range("H24").Value = 20000
$ws.Range("J24").Value = 20000
$ws.Range("L24").Value = 20000
$ws.Range("N24").Value = -20460
# Row 112
$ws.Range("H112").Value = 26675.4
$ws.Range("J112").Value = 26675.4
$ws.Range("L112").Value = 26675.4
$ws.Range("N112").Value = -29629.4
# Row 113
$ws.Range("H113").Value = 1086.75
$ws.Range("I113").Value = 933.9
$ws.Range("J113").Value = 1239.6
$ws.Range("K113").Value = 2801.7
$ws.Range("L113").Value = 3718.8
$ws.Range("M113").Value = -631.6999999999998
$ws.Range("N113").Value = -8058.799999999999
# Row 132
$ws.Range("H132").Value = 1561.4375
$ws.Range("I132").Value = 1373.6666
$ws.Range("K132").Value = 4120.9998
$ws.Range("M132").Value = -1590.9998
# Row 136
$ws.Range("H136").Value = 2210.4546
$ws.Range("I136").Value = 1818.3889
$ws.Range("J136").Value = 3974.75
$ws.Range("K136").Value = 5455.1667
$ws.Range("L136").Value = 11924.25
$ws.Range("M136").Value = -2905.1667
$ws.Range("N136").Value = -17024.25
